$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the database config values (graph alignment fix: bump version numbers)
$ws.Range("B5").Value = "highlands5"
$ws.Range("C5").Value = "survey2"

# Move the active selection up one row, from C7 to C6
$ws.Range("C6").Select()
